$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.125.75"
$ws.Range("E2").Value = "  +0.17%  "
$ws.Range("D3").Value = "1.874.88"
$ws.Range("E3").Value = "  +0.02%  "
$ws.Range("D4").Value = "'1.003"
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "'312.97"
$ws.Range("E5").Value = "  +0.02%  "
$ws.Range("E6").Value = "  -0.14%  "
$ws.Range("D7").Value = "'0.5134"
$ws.Range("E7").Value = "  +1.53%  "
$ws.Range("D8").Value = "'0.3892"
$ws.Range("E8").Value = "  +1.79%  "
$ws.Range("D9").Value = "'0.08393"
$ws.Range("E9").Value = "  -0.07%  "
$ws.Range("E10").Value = "  +0.37%  "
$ws.Range("D11").Value = "'41.56"
$ws.Range("E11").Value = "  -0.12%  "
$ws.Range("D12").Value = "'6.213"
$ws.Range("E12").Value = "  -0.15%  "
$ws.Range("D13").Value = "'20.73"
$ws.Range("E13").Value = "  +1.02%  "
$ws.Range("D14").Value = "1.868.12"
$ws.Range("E14").Value = "  -0.62%  "
$ws.Range("D15").Value = "'7.290"
$ws.Range("E15").Value = "  +1.31%  "
$ws.Range("D16").Value = "'1.003"
$ws.Range("E16").Value = "  -0.05%  "
$ws.Range("D17").Value = "'0.00001109"
$ws.Range("E17").Value = "  +1.05%  "
$ws.Range("D18").Value = "'90.98"
$ws.Range("E18").Value = "  -0.06%  "
$ws.Range("D19").Value = "'0.06650"
$ws.Range("E19").Value = "  -0.06%  "
$ws.Range("E20").Value = "  -1.60%  "
$ws.Range("D21").Value = "'1.002"
$ws.Range("E21").Value = "  -0.09%  "
$ws.Range("D22").Value = "'6.056"
$ws.Range("E22").Value = "  -0.10%  "
$ws.Range("D23").Value = "28.168.31"
$ws.Range("E23").Value = "  +0.20%  "
$ws.Range("D24").Value = "'11.14"
$ws.Range("E24").Value = "  -0.40%  "
$ws.Range("D25").Value = "'2.250"
$ws.Range("E25").Value = "  -0.61%  "
$ws.Range("D26").Value = "2.082.15"
$ws.Range("E26").Value = "  -0.69%  "
$ws.Range("D27").Value = "'2.501"
$ws.Range("E27").Value = "  -2.97%  "
$ws.Range("D28").Value = "'158.38"
$ws.Range("E28").Value = "  +0.71%  "
$ws.Range("D30").Value = "'125.08"
$ws.Range("E30").Value = "  -1.17%  "
$ws.Range("E31").Value = "  +1.17%  "
$ws.Range("D32").Value = "'1.041"
$ws.Range("E32").Value = "  -0.71%  "
$ws.Range("D33").Value = "'5.884"
$ws.Range("E33").Value = "  +4.69%  "
$ws.Range("D34").Value = "'3.600"
$ws.Range("E34").Value = "  -0.42%  "
$ws.Range("D35").Value = "'9.755"
$ws.Range("E35").Value = "  +0.53%  "
$ws.Range("D36").Value = "'0.02441"
$ws.Range("E36").Value = "  -0.55%  "
$ws.Range("D37").Value = "'0.06546"
$ws.Range("E37").Value = "  -0.05%  "
$ws.Range("D38").Value = "'0.2186"
$ws.Range("E38").Value = "  +0.82%  "
$ws.Range("D39").Value = "'1.210"
$ws.Range("E39").Value = "  -0.72%  "
$ws.Range("D40").Value = "'0.6514"
$ws.Range("E40").Value = "  +0.10%  "
$ws.Range("D41").Value = "'5.025"
$ws.Range("E41").Value = "  +2.69%  "
$ws.Range("D42").Value = "'1.227"
$ws.Range("E42").Value = "  -1.74%  "
$ws.Range("E43").Value = "  -0.03%  "
$ws.Range("D44").Value = "'0.6107"
$ws.Range("E44").Value = "  -1.32%  "
$ws.Range("D45").Value = "'13.12"
$ws.Range("E45").Value = "  -0.16%  "
$ws.Range("D46").Value = "'1.282"
$ws.Range("E46").Value = "  -1.57%  "
$ws.Range("D47").Value = "'3.673"
$ws.Range("E47").Value = "  -0.19%  "
$ws.Range("D48").Value = "'2.012"
$ws.Range("E48").Value = "  -0.05%  "
$ws.Range("E49").Value = "  -0.11%  "
$ws.Range("D50").Value = "'121.67"
$ws.Range("E50").Value = "  +0.57%  "
$ws.Range("D51").Value = "'77.95"
$ws.Range("E51").Value = "  -3.21%  "
